$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MigratoryModel_TableauData")
$ws.Activate()

# Rename Model values (column H) for the two AMNH language model blocks.
$ws.Range("H442:H471").Value = "AMNH - Language (sw-11-H)"
$ws.Range("H472:H501").Value = "AMNH - Language (sw-all5-H)"

# Update the current selection to reflect the last edited range.
$ws.Range("H472:H501").Select()
